# Apply updated mean_probability values (column C) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0.55
    3  = 0.55
    4  = 0.55
    5  = 0.55
    7  = 0.55
    8  = 0.55
    9  = 0.55
    10 = 0.55
    11 = 0.55
    12 = 0.55
    13 = 0.55
    14 = 0.55
    15 = 0.55
    16 = 0.55
    17 = 0.55
    18 = 0.55
    20 = 0.5071428571428571
    22 = 0.5333333333333333
    23 = 0.5333333333333333
    24 = 0.5333333333333333
    25 = 0.5333333333333333
    26 = 0.5333333333333333
    27 = 0.5333333333333333
    28 = 0.5333333333333333
    29 = 0.5333333333333333
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
